$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, pushing existing rows 12-47 down to 13-48.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = "2022-04-21"
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112022
$ws.Range("G12").Value = "Arveja Verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 23000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 24000
$ws.Range("N12").Value = "$/malla 25 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 960
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
